# Update v2 of the "Controle" sheet in the SIGLA data-processing workbook.
#
# Changes:
#  - B2: "9 x BACKUP 8GB RAM 8vCPU (WEBSERVER)"  -> "9 x BACKUP 16GB RAM 8vCPU (WEBSERVER)"
#  - B3: "4 x BACKUP 16GB RAM 16vCPU (WEBSERVER)" -> "4 x BACKUP 32GB RAM 16vCPU (WEBSERVER)"
#  - B5: "3 x BACKUP 24GB RAM 8vCPU (WEBSERVER)"  -> "3 x BACKUP 32GB RAM 8vCPU (WEBSERVER)"
#  - New (empty) cells materialised for columns that previously had no entry
#    on rows 2-5: H2, M2, N2, M3, N3, H4, M4, N4, H5, M5, N5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controle")

# --- text updates -----------------------------------------------------
$ws.Range("B2").Value = "9 x BACKUP 16GB RAM 8vCPU (WEBSERVER)"
$ws.Range("B3").Value = "4 x BACKUP 32GB RAM 16vCPU (WEBSERVER)"
$ws.Range("B5").Value = "3 x BACKUP 32GB RAM 8vCPU (WEBSERVER)"

# --- materialise empty string cells ------------------------------------
# Writing "" clears/omits the cell entirely, so instead we write a single
# apostrophe (Excel's "treat as text" prefix for an empty value) and then
# reset the style back to Normal so no stray formatting is left behind.
# The net effect is an empty-text cell that is actually present in the
# sheet, matching the target layout.
$emptyCells = @("H2", "M2", "N2", "M3", "N3", "H4", "M4", "N4", "H5", "M5", "N5")
foreach ($addr in $emptyCells) {
    $cell = $ws.Range($addr)
    $cell.Value = "'"
    $cell.Style = "Normal"
}
